# Auto-generated cell updates applying the cryptos.xlsx diff.
# NumberFormat is forced to Text ("@") before each assignment so that
# numeric-looking strings (prices like "353.98", "1.00", "0.0437", ...)
# are preserved verbatim as text instead of being auto-converted to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.039.24'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.004.67'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '353.98'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.78'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.20%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.00%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.04'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.83%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.65%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.98%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.478.80'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.61'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.999.44'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '52.123.90'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.13%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.54'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0971'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.10'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '263.68'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.56%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.48'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.47'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +6.98%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '35.99'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.28%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +15.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.10'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0437'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.84'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.82%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.58'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.23%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.29'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.119.88'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.75%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.245'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.77%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.904'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.65%  '
